$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 (Trey Murphy III is row 8, so this is the row right after it):
# "Royce O'Neale" -> "Ochai Agbaji", position "SF,PF" -> "SG,SF", team "Phoenix Suns" -> "Toronto Raptors"
$ws.Range("A9").Value = "Ochai Agbaji"
$ws.Range("B9").Value = "SG,SF"
$ws.Range("C9").Value = "Toronto Raptors"
